# Update the "Occurrence" column (E) values from the old list of numbers
# to the new list that has ", 85" appended at the end.
#
# In the original workbook every data row (E2:E13) shares the exact same
# string value (they all point at the same shared-string table entry), so
# we update every one of those cells to keep them consistent with the
# underlying shared string that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldValue = "19, 21, 23, 25, 27, 31, 33, 35, 39, 41, 45, 49, 53, 57, 61, 65, 73, 81"
$newValue = "19, 21, 23, 25, 27, 31, 33, 35, 39, 41, 45, 49, 53, 57, 61, 65, 73, 81, 85"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value()
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
